$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Rebuild the two hyperlinks up front: A2 keeps pointing at the server
# URL, and the A3 link grows to cover the new A4 row as well. The
# engine only supports a whole-sheet hyperlink wipe, so clear both and
# recreate them (A2 alone, A3:A4 merged).
# ------------------------------------------------------------------
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "http://172.191.4.85/TestCollection") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3:A4"), "http://172.191.4.85/TestCollection", "", "", "http://128.0.0.1/TestCollection") | Out-Null

# --- Row 2: Server URL / Project Name / PAT (Repository Name cleared) ---
$ws.Range("A2").Value = "http://128.0.0.1/TestCollection"
$ws.Range("B2").Value = "project1"
$ws.Range("C2").Value = "adad87adad8ds4449m434344mmnbnbb43434"
$ws.Range("D2").ClearContents()

# --- Row 3: Server URL / Project Name / PAT / Repository Name ---
$ws.Range("A3").Value = "http://128.0.0.1/TestCollection"
$ws.Range("B3").Value = "project2"
$ws.Range("C3").Value = "adad87adad8ds4449m434344mmnbnbb43434"
$ws.Range("D3").Value = "repo1"

# --- Row 4 (new): Server URL / Project Name / PAT / Repository Name ---
# Bring over C3's number/font formatting before writing C4's value.
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)

$ws.Range("A4").Value = "http://128.0.0.1/TestCollection"
$ws.Range("B4").Value = "project3"
$ws.Range("C4").Value = "adad87adad8ds4449m434344mmnbnbb43434"
$ws.Range("D4").Value = "repo2"

# Re-assert the Hyperlink cell style on A2:A4 - adding the hyperlinks
# above nudges the targeted cell onto a freshly duplicated style
# record; setting the named style again snaps it back onto the
# original shared "Hyperlink" xf so no new formatting is introduced.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("A4").Style = "Hyperlink"

Write-Output "done"
